$d = $word.ActiveDocument
$paras = $d.Paragraphs

# Locate the paragraph that ends the bibliography section.
$anchorIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "*Artigos de revistas especializadas*") {
        $anchorIndex = $i
        break
    }
}

# The three paragraphs right after it (a blank line, the "Ver no Jupiter..."
# line and the "(c) 2020 ..." footer line) are removed in their entirety,
# while the blank paragraph and page-break paragraph further below stay.
$startPara = $paras.Item($anchorIndex + 1)
$endPara = $paras.Item($anchorIndex + 3)

$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()
